$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(9, 2).Value = "51cd7587f6b732c6eff305629aeadfd8"  # 05-050305TC
$ws.Cells.Item(11, 2).Value = "f2d8e4b51f987ad30c3ce9202f61284c"  # 05-050301A
$ws.Cells.Item(15, 2).Value = "ff33e05ce4124f5cbf53bf885c4af68d"  # 05-050207TP
$ws.Cells.Item(17, 2).Value = "be08aba23e43185b17caf928993046e4"  # 05-050305TP
$ws.Cells.Item(24, 2).Value = "5f0f80dc3ac770be07bf527bf4b48261"  # 05-050316TC
$ws.Cells.Item(29, 2).Value = "763c606937babd4da942d51c5671be4a"  # 05-050302A
$ws.Cells.Item(34, 2).Value = "37356ee34e7335b4597df2382e2a9bb3"  # 05-050316TP
$ws.Cells.Item(89, 2).Value = "9becbb1e4390650867e2a1d98dc8fb05"  # 05-050104A
$ws.Cells.Item(99, 2).Value = "43b88bf3307da011bc0998dbad5276a5"  # 05-050101A
$ws.Cells.Item(110, 2).Value = "619b4c12b3d41217370712afe12beea7"  # 05-050102A
$ws.Cells.Item(121, 2).Value = "1f698a3405bd15320c1955fb18b791a3"  # 05-050301TP
$ws.Cells.Item(133, 2).Value = "c8fe16b82802716bc3cf7fca3460657c"  # 05-050312TP
$ws.Cells.Item(136, 2).Value = "185bea02a4e2c5ff9756796cf930cb7a"  # 05-050312TC
$ws.Cells.Item(160, 2).Value = "30927c55b50630f7226ca8614243884f"  # 05-050203TP
$ws.Cells.Item(163, 2).Value = "3bbcea997c54685cddb80bc361e32ff6"  # 05-050308A
$ws.Cells.Item(170, 2).Value = "c566b9b1d1f66fa3fd052765cc85f92d"  # 05-050203TC
$ws.Cells.Item(176, 2).Value = "856a7a25028d800764ae2dc8a4dd3c5d"  # 05-050303TP
$ws.Cells.Item(181, 2).Value = "0aeceaedd3b4392d28b6ced93687f649"  # 05-050303TC
$ws.Cells.Item(184, 2).Value = "594effaede05e0c92cfc818960a7dae9"  # 05-050305A
$ws.Cells.Item(192, 2).Value = "ac2b19875e22bcea34db82975626dda9"  # 05-050314TP
$ws.Cells.Item(199, 2).Value = "638207d29cecb4dab9e11f72d0d9906c"  # 05-050314TC
$ws.Cells.Item(201, 2).Value = "d616a6a4bc33792ef4f73a2c8b60bb51"  # 05-050306A
$ws.Cells.Item(214, 2).Value = "3350a2f269f5aa9197ab70bfd89b98ac"  # 05-050303A
$ws.Cells.Item(229, 2).Value = "f7dca3d50595a80148e0994682fe50ca"  # 05-050205TP
$ws.Cells.Item(230, 2).Value = "3ba9785b382422979b88e10bcbb1e127"  # 05-050304A
$ws.Cells.Item(234, 2).Value = "f442fd8e6eb755b52e541f5509e1a89f"  # 05-050205TC
$ws.Cells.Item(284, 2).Value = "45865d328aa4c09e7908e6b683036d88"  # 05-050101TP
$ws.Cells.Item(345, 2).Value = "38a2e984b4e6e3225db6d45a4ca2e5b5"  # 05-050201TP
$ws.Cells.Item(470, 2).Value = "d0d6e1185552166477780d4dec26a7d2"  # 05-050204A
$ws.Cells.Item(471, 2).Value = "afbbc41902af0f586a4df26aa930d626"  # 05-050313A
$ws.Cells.Item(489, 2).Value = "34eb4fba8c47b138d9f1b2aebddb44ac"  # 05-050205A
$ws.Cells.Item(491, 2).Value = "8a26d1e6fea5ad546cc92ee5ececa97b"  # 05-050314A
$ws.Cells.Item(505, 2).Value = "581b8d3e661eb4b3e04ddfb924f5ec62"  # 05-050208TC
$ws.Cells.Item(514, 2).Value = "69e14cc3f1fcba16cedbb91771297a1d"  # 05-050311A
$ws.Cells.Item(520, 2).Value = "60ba7d73f034edb735dd980c736ccc77"  # 05-050306TP
$ws.Cells.Item(528, 2).Value = "b296ddf41fecd4d8db364b216f92aee6"  # 05-050317TC
$ws.Cells.Item(531, 2).Value = "6b576b5fe60fc372b35c4df3c857f680"  # 05-050203A
$ws.Cells.Item(539, 2).Value = "b5d5aeacc6075482da99345d8dfa9f1d"  # 05-050317TP
$ws.Cells.Item(563, 2).Value = "732acee5513bf6128fe232b1d3984b69"  # 05-050201A
$ws.Cells.Item(579, 2).Value = "8b3033894d2c254c5ce25ddcc8d9e1d7"  # 05-050308TC
$ws.Cells.Item(588, 2).Value = "aa6129345b89ae214245282954d53b14"  # 05-050308TP
$ws.Cells.Item(632, 2).Value = "9711971db29b5bcd80a547b33ef49502"  # 05-050204TP
$ws.Cells.Item(643, 2).Value = "5ef1f634c309318a4d5c86122cdea80c"  # 05-050204TC
$ws.Cells.Item(645, 2).Value = "40a3ddd5a698e1ec278c0dfe84da4eb3"  # 05-050302TP
$ws.Cells.Item(665, 2).Value = "f72143de0917f37cd2702ec3118a6598"  # 05-050313TP
$ws.Cells.Item(671, 2).Value = "cd26e29331dc2588008e0b239924e1b4"  # 05-050313TC
$ws.Cells.Item(682, 2).Value = "4361a799ad462dd645713e10dacae2b5"  # 05-050317A
$ws.Cells.Item(696, 2).Value = "22a00d96591da67b31f3bb98464bad65"  # 05-050206TP
$ws.Cells.Item(701, 2).Value = "765303b626f163da6d61cd38f37360ed"  # 05-050206TC
$ws.Cells.Item(716, 2).Value = "841e1fb15e16a1c76a4e96830a33fbe4"  # 05-050304TC
$ws.Cells.Item(719, 2).Value = "462e1f0385f2fc9212acf0fa97e73d6d"  # 05-050206A
$ws.Cells.Item(720, 2).Value = "c09abe3a701ea1fce4061d837f4755e1"  # 05-050315A
$ws.Cells.Item(731, 2).Value = "caf8ca558cf92385e4a7602a55b21e53"  # 05-050304TP
$ws.Cells.Item(745, 2).Value = "e823c65b96cde202b8317746781c48a3"  # 05-050316A
$ws.Cells.Item(749, 2).Value = "74c6706958fb2ca6a70191ae5cb0cf45"  # 05-050207A
$ws.Cells.Item(758, 2).Value = "4d593788b26eb4705d5e2cf0f80e57bb"  # 05-050315TP
$ws.Cells.Item(785, 2).Value = "a126afa1c92f45988645e48cb0eb7123"  # 05-050102TP
$ws.Cells.Item(844, 2).Value = "d3cee495caa36d6a4f400c78c22f835f"  # 05-050104TP
$ws.Cells.Item(845, 2).Value = "38a3977ed5e0681da1529ab7c218b3de"  # 05-050202TC
$ws.Cells.Item(848, 2).Value = "31e56c1dac66cb42ec103546baebcc92"  # 05-050311TC
$ws.Cells.Item(853, 2).Value = "72c27219de084da91a0e32afbd0f9dde"  # 05-050311TP
$ws.Cells.Item(872, 2).Value = "9dfdade307c3989cf6e4dbf168845da0"  # 05-050309TC
